# New release build of the workbook: the "build on" timestamp embedded in
# several descriptive text cells is bumped from the January 30 build to the
# February 02 build (time also changes, EST in both cases).

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

# --- "About" sheet -----------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value2 = $about.Range("A2").Value2.Replace($oldStamp, $newStamp)
$about.Range("A6").Value2 = $about.Range("A6").Value2.Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet -----------------------------
$data = $wb.Worksheets.Item("Boundaries and methane sources")

$lastRow = $data.UsedRange.Rows.Count
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $data.Cells.Item($row, 19)  # column S = build_version
    $current = $cell.Value2
    if ($current -ne $null -and $current -is [string] -and $current.Contains($oldStamp)) {
        $cell.Value2 = $current.Replace($oldStamp, $newStamp)
    }
}
